# budget_upload_test.xlsx - "First set of budget test working."
#
# 1. Sheet "Budgets": cell A3 value 109189 -> 888888
# 2. Sheet "Budgets": selection moves from G6 to A3
# 3. Sheet "Budgets": page setup gets an explicit paper size / orientation
#    (A4 / portrait -> paperSize=9, orientation=portrait)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cell edit ---------------------------------------------------------
$ws.Range("A3").Value = 888888

# --- page setup ----------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- selection (do this last so it "sticks" as the active selection) ---
$ws.Range("A3").Select() | Out-Null
